# Updated symbol list on Sat Feb  4 23:53:45 UTC 2023 with GitHub Actions
# Refresh of Price (column D) and Volume(1h) (column E) for the crypto ticker rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> updated Price / Volume(1h) text values.
# A leading apostrophe forces Excel to store the value as literal text
# (these columns hold formatted strings like "332.10" / "0.02%", not numbers).
$updates = @{
    2 = @{ D="332.10"; E="0.02%" }
    3 = @{ D="41.22"; E="0.13%" }
    4 = @{ D="5.693"; E="-0.71%" }
    5 = @{ D="0.08393"; E="2.34%" }
    6 = @{ D="8.810"; E="0.73%" }
    7 = @{ D="4.522"; E="0.24%" }
    8 = @{ E="-2.90%" }
    9 = @{ D="2.932"; E="-2.03%" }
    10 = @{ D="0.9245"; E="0.31%" }
    11 = @{ D="0.1245"; E="-0.10%" }
    12 = @{ D="0.1956"; E="0.11%" }
    13 = @{ D="0.09354"; E="-0.98%" }
    14 = @{ D="0.03969"; E="9.22%" }
    15 = @{ E="0.79%" }
    16 = @{ D="0.001299"; E="-0.14%" }
    17 = @{ D="0.006122"; E="-1.47%" }
    18 = @{ E="1.46%" }
    20 = @{ D="9.120"; E="9.70%" }
    21 = @{ D="0.1378"; E="-2.69%" }
    22 = @{ D="0.2641"; E="-0.36%" }
    23 = @{ D="0.04422"; E="-0.13%" }
    24 = @{ D="0.001247"; E="-1.02%" }
    25 = @{ D="0.004351"; E="0.82%" }
    26 = @{ D="0.0001195"; E="-3.68%" }
    27 = @{ D="0.0004013"; E="0.50%" }
    39 = @{ D="0.02797"; E="0.80%" }
    40 = @{ D="0.05526"; E="0.20%" }
    41 = @{ D="0.007940"; E="4.31%" }
    42 = @{ D="0.1429"; E="0.41%" }
    43 = @{ D="0.009010" }
    44 = @{ D="0.002177"; E="2.19%" }
    45 = @{ D="0.01011"; E="-14.61%" }
    46 = @{ D="0.00007175"; E="6.23%" }
    47 = @{ E="0.46%" }
    48 = @{ D="0.003466"; E="18.31%" }
    49 = @{ D="0.002291"; E="0.52%" }
    50 = @{ E="0.46%" }
    51 = @{ E="0.46%" }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    if ($cells.ContainsKey("D")) {
        $ws.Range("D$row").Value = "'" + $cells["D"]
    }
    if ($cells.ContainsKey("E")) {
        $ws.Range("E$row").Value = "'" + $cells["E"]
    }
}
